$d = $word.ActiveDocument

# 1) Paragraph 1: "BBBBBBBBBBBBBB" -> "Person 1"
$d.Content.Find.Execute("BBBBBBBBBBBBBB", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Person 1", 2)

# 2) Remove the empty paragraph (old paragraph 2) and the old paragraph 3
#    ("CCCCCCCCCCCCCCCCCCCCCcccccc", wrapped in proofErr spell-check markers)
#    entirely -- deleting the whole paragraph range (including its proofErr
#    runs and paragraph mark) is what drops the proofErr elements.
$r = $d.Range($d.Paragraphs(2).Range.Start, $d.Paragraphs(3).Range.End)
$r.Delete()

# 3) The trailing paragraph left behind becomes "Person 2 " (note trailing space)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Text = "Person 2 "

# 4) Append "Person 3"
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Text = "Person 3"

# 5) Append "Person 4"
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Text = "Person 4"
